# Add "Number of siblings" and "Race" columns to the Data sheet, and their
# corresponding rows in the Codebook sheet.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsCodebook = $wb.Worksheets.Item("Codebook")

# ---- Data sheet: new headers ----
$wsData.Range("D1").Value = "Number of siblings"
$wsData.Range("E1").Value = "Race"
$wsData.Range("D1:E1").Font.Bold = $true

# ---- Data sheet: new column values ----
$siblings = @(2, 3, 4, 1, 2, 0, 2, 3, 1, 3, 3, 1, 2, 4)
$race = @("White", "Black", "Hispanic", "Asian", "Asian", "Other", "Hispanic", "Other", "Asian", "White", "White", "Black", "Black", "Hispanic")

for ($i = 0; $i -lt $siblings.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 4).Value = $siblings[$i]
    $wsData.Cells.Item($row, 5).Value = $race[$i]
}

# ---- Codebook sheet: new rows describing the new variables ----
$wsCodebook.Range("A5").Value = "Number of siblings"
$wsCodebook.Range("B5").Value = "Integer number for number of siblings"
$wsCodebook.Range("C5").Value = "numeric value >=0"

$wsCodebook.Range("A6").Value = "Race"
$wsCodebook.Range("B6").Value = "Self-identification of race: White and Black refer to Non-Hispanic White and Non-Hispanic Black, Hispanic refers to both Hispanic and Latino individuals, Asians refer to South and East Asians; those that do not identify as any of these are classified as Other"
$wsCodebook.Range("C6").Value = "White/Black/Hispanic/Asian/Other"

# ---- Selections / active sheet ----
$wsData.Range("F5").Select()
$wsCodebook.Activate()
$wsCodebook.Range("A5:C6").Select()
